$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.899.38"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.861.30"
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'244.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'0.6590"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'48.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.2998"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'24.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.12%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07643"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.864.78"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.082"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.6926"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'83.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000009727"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'6.146"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.43%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "29.915.32"
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.121.43"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'237.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'12.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'7.829"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.48%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'1.001"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1446"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'158.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'8.621"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'17.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'1.497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.06053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.291"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.171"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.097"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.877"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.184"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.7360"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.609"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.209.24"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.9173"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.035.22"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'67.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.20%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'101.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.95%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000122"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4081"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.234"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "

Write-Host "Updated cryptos list"
